$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1007048317"
$ws.Range("D16").Value = "MARCO POLO HERRERA BERMEJO"
$ws.Range("E16").Value = "2210"
$ws.Range("F16").Value = 60000
$ws.Range("G16").Value = 1500000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143390341"
$ws.Range("D17").Value = "DANIEL CAMILO CASTRO REALES"
$ws.Range("E17").Value = "2210"
$ws.Range("F17").Value = 60000
$ws.Range("G17").Value = 1500000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1007048317"
$ws.Range("D18").Value = "MARCO POLO HERRERA BERMEJO"
$ws.Range("E18").Value = "2211"
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1500000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143390341"
$ws.Range("D19").Value = "DANIEL CAMILO CASTRO REALES"
$ws.Range("E19").Value = "2211"
$ws.Range("F19").Value = 60000
$ws.Range("G19").Value = 1500000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1007048317"
$ws.Range("D20").Value = "MARCO POLO HERRERA BERMEJO"
$ws.Range("E20").Value = "2212"
$ws.Range("F20").Value = 60000
$ws.Range("G20").Value = 1500000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143390341"
$ws.Range("D21").Value = "DANIEL CAMILO CASTRO REALES"
$ws.Range("E21").Value = "2212"
$ws.Range("F21").Value = 60000
$ws.Range("G21").Value = 1500000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1007048317"
$ws.Range("D22").Value = "MARCO POLO HERRERA BERMEJO"
$ws.Range("E22").Value = "2301"
$ws.Range("F22").Value = 60000
$ws.Range("G22").Value = 1500000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143390341"
$ws.Range("D23").Value = "DANIEL CAMILO CASTRO REALES"
$ws.Range("E23").Value = "2301"
$ws.Range("F23").Value = 60000
$ws.Range("G23").Value = 1500000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1007048317"
$ws.Range("D24").Value = "MARCO POLO HERRERA BERMEJO"
$ws.Range("E24").Value = "2302"
$ws.Range("F24").Value = 60000
$ws.Range("G24").Value = 1500000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1143390341"
$ws.Range("D25").Value = "DANIEL CAMILO CASTRO REALES"
$ws.Range("E25").Value = "2302"
$ws.Range("F25").Value = 60000
$ws.Range("G25").Value = 1500000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1007048317"
$ws.Range("D26").Value = "MARCO POLO HERRERA BERMEJO"
$ws.Range("E26").Value = "2303"
$ws.Range("F26").Value = 60000
$ws.Range("G26").Value = 1500000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1143390341"
$ws.Range("D27").Value = "DANIEL CAMILO CASTRO REALES"
$ws.Range("E27").Value = "2303"
$ws.Range("F27").Value = 60000
$ws.Range("G27").Value = 1500000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1007048317"
$ws.Range("D28").Value = "MARCO POLO HERRERA BERMEJO"
$ws.Range("E28").Value = "2304"
$ws.Range("F28").Value = 54000
$ws.Range("G28").Value = 1500000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1143390341"
$ws.Range("D29").Value = "DANIEL CAMILO CASTRO REALES"
$ws.Range("E29").Value = "2304"
$ws.Range("F29").Value = 54000
$ws.Range("G29").Value = 1500000

